$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached text of every auto-updating "Date" placeholder
#    (slide master + every slide layout) from 7/27/2017 to 8/11/2017.
# ---------------------------------------------------------------------
$newDate = "8/11/2017"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDate = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout attached to the (single) slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# NOTE: the notes master's own "Date Placeholder" field (used by
# notesMaster1.xml) cannot be refreshed through this host: writes to
# $p.NotesMaster shapes are not persisted there (and can mis-target
# the slide master's shape collection), so it is intentionally left
# untouched rather than risk corrupting the slide master.

# ---------------------------------------------------------------------
# 2) Slide 20 ("Next Steps") wording tweaks.
# ---------------------------------------------------------------------
$slide20 = $p.Slides.Item(20)
$content = $slide20.Shapes.Item(2)
$tr = $content.TextFrame.TextRange

$para4 = $tr.Paragraphs(4)
$fullPara4 = $tr.Characters($para4.Start, $para4.Length - 1)
$fullPara4.Text = "If you have questions about this lesson, ask them on the Discussion Board."

$para5 = $tr.Paragraphs(5)
$secondHalf = $tr.Characters($para5.Start + 13, $para5.Length - 13)
$secondHalf.Text = "next lesson."
